# Apply the described changes:
# 1. Rename the "congenital" sheet to "misc_long_term" (tab name + header cell A1).
# 2. On the "mental" sheet, remove the "Ocular" row, shifting the rows below it up.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "congenital" sheet ---
$congenital = $wb.Worksheets.Item("congenital")
$congenital.Name = "misc_long_term"
$congenital.Range("A1").Value = "misc_long_term"

# --- 2. Remove the "Ocular" row from the "mental" sheet ---
$mental = $wb.Worksheets.Item("mental")
$mental.Rows.Item(5).Delete()
